# Populate the SEC_Comm / SEC_Processes / MIN_IMP sheets with the new
# SOLAR / WIND_ON / WIND_OFF / NAT_GAS / ELC commodities + processes,
# and wire up the MIN_IMP lookup table that reads them back.

$wb = $excel.ActiveWorkbook

$wsComm  = $wb.Worksheets.Item("SEC_Comm")
$wsProc  = $wb.Worksheets.Item("SEC_Processes")
$wsMin   = $wb.Worksheets.Item("MIN_IMP")

# ---------------------------------------------------------------------
# SEC_Comm!B7:I11  -- "Define Commodities" table
# Columns: B=Csets, C=CommName, D=CommDesc, E=Unit, F=LimType,
#          G=CTSLvl, H=PeakTS, I=Ctype
# ---------------------------------------------------------------------
$wsComm.Range("B7").Value  = "NRG"
$wsComm.Range("C7").Value  = "SOLAR"
$wsComm.Range("D7").Value  = "Solar irradiance process"
$wsComm.Range("E7").Value  = "PJ"
$wsComm.Range("G7").Value  = "DAYNITE"

$wsComm.Range("B8").Value  = "NRG"
$wsComm.Range("C8").Value  = "WIND_ON"
$wsComm.Range("D8").Value  = "Wind onshore"
$wsComm.Range("E8").Value  = "PJ"
$wsComm.Range("G8").Value  = "DAYNITE"

$wsComm.Range("B9").Value  = "NRG"
$wsComm.Range("C9").Value  = "WIND_OFF"
$wsComm.Range("D9").Value  = "Wind offshore"
$wsComm.Range("E9").Value  = "PJ"
$wsComm.Range("G9").Value  = "DAYNITE"

$wsComm.Range("B10").Value = "NRG"
$wsComm.Range("C10").Value = "NAT_GAS"
$wsComm.Range("D10").Value = "Nat gas"
$wsComm.Range("E10").Value = "PJ"

$wsComm.Range("B11").Value = "NRG"
$wsComm.Range("C11").Value = "ELC"
$wsComm.Range("D11").Value = "Electricity"
$wsComm.Range("E11").Value = "PJ"
$wsComm.Range("G11").Value = "DAYNITE"
$wsComm.Range("I11").Value = "ELC"

# ---------------------------------------------------------------------
# SEC_Processes!B7:H11  -- "Define Processes" table
# Columns: B=Sets, D=TechName, F=Tact, G=Tcap, H=Tslvl
# ---------------------------------------------------------------------
$wsProc.Range("B7").Value  = "MIN"
$wsProc.Range("D7").Value  = "MIN_SOLAR"
$wsProc.Range("F7").Value  = "PJ"
$wsProc.Range("G7").Value  = "PJa"
$wsProc.Range("H7").Value  = "DAYNITE"

$wsProc.Range("B8").Value  = "MIN"
$wsProc.Range("D8").Value  = "MIN_WIND_ON"
$wsProc.Range("F8").Value  = "PJ"
$wsProc.Range("G8").Value  = "PJa"
$wsProc.Range("H8").Value  = "DAYNITE"

$wsProc.Range("B9").Value  = "MIN"
$wsProc.Range("D9").Value  = "MIN_WIND_OFF"
$wsProc.Range("F9").Value  = "PJ"
$wsProc.Range("G9").Value  = "PJa"
$wsProc.Range("H9").Value  = "DAYNITE"

$wsProc.Range("B10").Value = "IMP"
$wsProc.Range("D10").Value = "IMP_NAT_GAS"
$wsProc.Range("F10").Value = "PJ"
$wsProc.Range("G10").Value = "PJa"

$wsProc.Range("B11").Value = "IMP"
$wsProc.Range("D11").Value = "IMP_ELC"
$wsProc.Range("F11").Value = "PJ"
$wsProc.Range("G11").Value = "PJa"
$wsProc.Range("H11").Value = "DAYNITE"

# ---------------------------------------------------------------------
# MIN_IMP!E5:F5  -- extra headers
# ---------------------------------------------------------------------
$wsMin.Range("E5").Value = "COST"
$wsMin.Range("F5").Value = "CUM"

# MIN_IMP!B7:F11  -- lookups back into SEC_Processes / SEC_Comm, plus
# the extraction-cost / cumulative-bound values.
# D7:D11 is one dynamic/legacy array formula spilling SEC_Comm!C7:C11,
# so set it once over the whole range (fills D7..D11 in one shot).
$wsMin.Range("D7:D11").FormulaArray = "=SEC_Comm!C7:C11"

$wsMin.Range("B7").Formula = "=SEC_Processes!D7"
$wsMin.Range("E7").Value = 0.00001

$wsMin.Range("B8").Formula = "=SEC_Processes!D8"
$wsMin.Range("E8").Value = 0.00001

$wsMin.Range("B9").Formula = "=SEC_Processes!D9"
$wsMin.Range("E9").Value = 0.00001

$wsMin.Range("B10").Formula = "=SEC_Processes!D10"
$wsMin.Range("E10").Value = 10
$wsMin.Range("I10").Value = "We put 0.0001 as a ""very small"" value that is not 0"

$wsMin.Range("B11").Formula = "=SEC_Processes!D11"
$wsMin.Range("E11").Value = 50
